# Helper: set a cell to a value; values prefixed with "TXT:" are forced to
# text (so numeric-looking strings like fund codes "006113" or percentages
# "85.14" are kept as text instead of being auto-parsed into numbers), while
# plain values are written as-is (numbers stay numbers).
function Set-CellVal($cell, $val) {
    if ($val -is [string] -and $val.StartsWith("TXT:")) {
        $text = $val.Substring(4)
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet (sheet1): insert a new "2022-Q3" row right after the
#    header, shifting the existing quarterly summary rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

for ($r = 8; $r -ge 2; $r--) {
    $destRow = $r + 1
    $total.Cells.Item($destRow, 1).Value = $total.Cells.Item($r, 1).Value()
    $total.Cells.Item($destRow, 2).Value = $total.Cells.Item($r, 2).Value()
    $total.Cells.Item($destRow, 3).Value = $total.Cells.Item($r, 3).Value()
    $total.Cells.Item($destRow, 4).Value = $total.Cells.Item($r, 4).Value()
}
# Row 9 is brand new -- give column A the same style as the rest of the
# column (centered bold style used throughout column A of this sheet).
$total.Cells.Item(8, 1).Copy()
$total.Cells.Item(9, 1).PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 20
$total.Cells.Item(2, 4).Value = 9

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q3" fund-holdings sheet. Duplicate the existing
#    "2022-Q2" sheet (same column layout/formatting) and place the copy
#    right before it, then overwrite its data.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$rows3 = @(
    @(0, 'TXT:006113', 'TXT:汇添富创新医药主题混合', 'TXT:85.14', 'TXT:75.99', 'TXT:3.77', 'TXT:3.2098', 8),
    @(1, 'TXT:470006', 'TXT:汇添富医药保健混合A', 'TXT:46.47', 'TXT:84.46', 'TXT:4.24', 'TXT:1.9703', 4),
    @(2, 'TXT:009664', 'TXT:汇添富医疗积极成长一年持有期混合A', 'TXT:29.13', 'TXT:66.68', 'TXT:5.00', 'TXT:1.4565', 5),
    @(3, 'TXT:015115', 'TXT:汇添富中国高端制造股票D', 'TXT:24.58', 'TXT:83.58', 'TXT:3.31', 'TXT:0.8136', 9),
    @(4, 'TXT:001725', 'TXT:汇添富中国高端制造股票A', 'TXT:24.16', 'TXT:83.58', 'TXT:3.31', 'TXT:0.7997', 9),
    @(5, 'TXT:007639', 'TXT:汇添富竞争优势灵活配置混合', 'TXT:7.29', 'TXT:86.04', 'TXT:4.59', 'TXT:0.3346', 5),
    @(6, 'TXT:009665', 'TXT:汇添富医疗积极成长一年持有期混合C', 'TXT:3.63', 'TXT:66.68', 'TXT:5.00', 'TXT:0.1815', 5),
    @(7, 'TXT:014126', 'TXT:华夏中证1000指数增强C', 'TXT:8.78', 'TXT:89.62', 'TXT:0.82', 'TXT:0.0720', 5),
    @(8, 'TXT:501063', 'TXT:汇添富悦享定期开放混合', 'TXT:2.19', 'TXT:66.28', 'TXT:3.06', 'TXT:0.0670', 10),
    @(9, 'TXT:515760', 'TXT:华夏中证浙江国资创新发展ETF', 'TXT:2.04', 'TXT:99.57', 'TXT:2.90', 'TXT:0.0592', 10),
    @(10, 'TXT:012430', 'TXT:农银汇理瑞康6个月持有期混合', 'TXT:1.30', 'TXT:24.44', 'TXT:1.10', 'TXT:0.0143', 5),
    @(11, 'TXT:014125', 'TXT:华夏中证1000指数增强A', 'TXT:0.97', 'TXT:89.62', 'TXT:0.82', 'TXT:0.0080', 5),
    @(12, 'TXT:015466', 'TXT:太平中证1000指数增强A', 'TXT:0.37', 'TXT:92.23', 'TXT:1.06', 'TXT:0.0039', 5),
    @(13, 'TXT:014820', 'TXT:华安创新医药锐选量化股票A', 'TXT:0.10', 'TXT:91.42', 'TXT:2.35', 'TXT:0.0024', 9),
    @(14, 'TXT:015114', 'TXT:汇添富中国高端制造股票C', 'TXT:0.04', 'TXT:83.58', 'TXT:3.31', 'TXT:0.0013', 9),
    @(15, 'TXT:006143', 'TXT:恒生前海中证质量成长低波动指数A', 'TXT:0.05', 'TXT:93.33', 'TXT:2.23', 'TXT:0.0011', 10),
    @(16, 'TXT:014821', 'TXT:华安创新医药锐选量化股票C', 'TXT:0.03', 'TXT:91.42', 'TXT:2.35', 'TXT:0.0007', 9),
    @(17, 'TXT:015467', 'TXT:太平中证1000指数增强C', 'TXT:0.02', 'TXT:92.23', 'TXT:1.06', 'TXT:0.0002', 5),
    @(18, 'TXT:006144', 'TXT:恒生前海中证质量成长低波动指数C', 'TXT:0.01', 'TXT:93.33', 'TXT:2.23', 'TXT:0.0002', 10),
    @(19, 'TXT:960015', 'TXT:汇添富医药保健混合O', 'TXT:0.00', 'TXT:84.46', 'TXT:4.24', '0', 4)
)

$rowIdx = 2
foreach ($row in $rows3) {
    Set-CellVal $q3.Cells.Item($rowIdx, 1) $row[0]
    Set-CellVal $q3.Cells.Item($rowIdx, 2) $row[1]
    Set-CellVal $q3.Cells.Item($rowIdx, 3) $row[2]
    Set-CellVal $q3.Cells.Item($rowIdx, 4) $row[3]
    Set-CellVal $q3.Cells.Item($rowIdx, 5) $row[4]
    Set-CellVal $q3.Cells.Item($rowIdx, 6) $row[5]
    Set-CellVal $q3.Cells.Item($rowIdx, 7) $row[6]
    Set-CellVal $q3.Cells.Item($rowIdx, 8) $row[7]
    $rowIdx = $rowIdx + 1
}

# Rows 12-21 are new (the source sheet only had 11 rows); give column A
# the same style ("s=2" style used for the rest of column A) as row 11.
$q3.Cells.Item(11, 1).Copy()
$q3.Range("A12:A21").PasteSpecial(-4122)

# Re-apply the A-column values after the format paste (PasteSpecial with
# formats only should not disturb values, but set again to be safe).
$rowIdx = 2
foreach ($row in $rows3) {
    $q3.Cells.Item($rowIdx, 1).Value = $row[0]
    $rowIdx = $rowIdx + 1
}

Write-Host "Edit complete"
